# Auto-generated Excel COM-interop script to apply Chocobo_Profits.xlsx numeric updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1669.9487
$ws.Range("I112").Value = 941.5
$ws.Range("J112").Value = 1802.3939
$ws.Range("K112").Value = 2824.5
$ws.Range("L112").Value = 5407.1817
$ws.Range("M112").Value = -1716.5
$ws.Range("N112").Value = -7623.1817

$ws.Range("H132").Value = 346105.6
$ws.Range("I132").Value = 193088.72
$ws.Range("J132").Value = 1672251.9
$ws.Range("K132").Value = 579266.16
$ws.Range("L132").Value = 5016755.699999999
$ws.Range("M132").Value = -576736.16
$ws.Range("N132").Value = -5021815.699999999

$ws.Range("H137").Value = 2999.639
$ws.Range("I137").Value = 1504.2727
$ws.Range("J137").Value = 5349.5
$ws.Range("K137").Value = 4512.8181
$ws.Range("L137").Value = 16048.5
$ws.Range("M137").Value = -1962.8181
$ws.Range("N137").Value = -21148.5

$ws.Range("H138").Value = 3199.93
$ws.Range("I138").Value = 707.7222
$ws.Range("J138").Value = 4601.797
$ws.Range("K138").Value = 2123.1666
$ws.Range("L138").Value = 13805.391
$ws.Range("M138").Value = 3016.8334
$ws.Range("N138").Value = -24085.391

$ws.Range("H141").Value = 4516.8945
$ws.Range("I141").Value = 4586.352
$ws.Range("J141").Value = 3266.6667
$ws.Range("K141").Value = 13759.056
$ws.Range("L141").Value = 9800.000100000001
$ws.Range("M141").Value = -8579.056
$ws.Range("N141").Value = -20160.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5256.385
$ws.Range("I32").Value = 4129.06
$ws.Range("K32").Value = 4129.06
$ws.Range("M32").Value = -3842.06

$ws.Range("H61").Value = 980.7561
$ws.Range("I61").Value = 798.14703
$ws.Range("J61").Value = 1867.7142
$ws.Range("K61").Value = 798.14703
$ws.Range("L61").Value = 1867.7142
$ws.Range("M61").Value = -586.14703
$ws.Range("N61").Value = -2291.7142

$ws.Range("H122").Value = 1972.08
$ws.Range("I122").Value = 1189.6
$ws.Range("K122").Value = 3568.8
$ws.Range("M122").Value = -1118.8

$ws.Range("H132").Value = 2206.1353
$ws.Range("I132").Value = 1208.5927
$ws.Range("J132").Value = 4899.5
$ws.Range("K132").Value = 3625.7781
$ws.Range("L132").Value = 14698.5
$ws.Range("M132").Value = -1095.7781
$ws.Range("N132").Value = -19758.5

$ws.Range("H136").Value = 980.7561
$ws.Range("I136").Value = 798.14703
$ws.Range("J136").Value = 1867.7142
$ws.Range("K136").Value = 2394.44109
$ws.Range("L136").Value = 5603.142599999999
$ws.Range("M136").Value = 155.5589100000002
$ws.Range("N136").Value = -10703.1426

$ws.Range("H137").Value = 41780
$ws.Range("J137").Value = 41780
$ws.Range("L137").Value = 41780
$ws.Range("N137").Value = -51980

$ws.Range("H141").Value = 48122.5
$ws.Range("J141").Value = 48122.5
$ws.Range("L141").Value = 48122.5
$ws.Range("N141").Value = -58482.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1412.2909
$ws.Range("I134").Value = 1013.7959
$ws.Range("K134").Value = 3041.3877
$ws.Range("M134").Value = -506.3876999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9617433
$ws.Range("I31").Value = 1134.7567
$ws.Range("K31").Value = 1134.7567
$ws.Range("M31").Value = -839.7566999999999

$ws.Range("H34").Value = 9617433
$ws.Range("I34").Value = 1134.7567
$ws.Range("K34").Value = 1134.7567
$ws.Range("M34").Value = -932.7566999999999

$ws.Range("H132").Value = 1175.8871
$ws.Range("I132").Value = 680.1111
$ws.Range("K132").Value = 2040.3333
$ws.Range("M132").Value = 489.6667000000002

$ws.Range("H134").Value = 1618.0588
$ws.Range("I134").Value = 700.56
$ws.Range("J134").Value = 4166.6665
$ws.Range("K134").Value = 2101.68
$ws.Range("L134").Value = 12499.9995
$ws.Range("M134").Value = 433.3200000000002
$ws.Range("N134").Value = -17569.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3288.2144
$ws.Range("I63").Value = 2603.6
$ws.Range("J63").Value = 4999.75
$ws.Range("K63").Value = 7810.799999999999
$ws.Range("L63").Value = 14999.25
$ws.Range("M63").Value = -7061.799999999999
$ws.Range("N63").Value = -16497.25

$ws.Range("H66").Value = 3288.2144
$ws.Range("I66").Value = 2603.6
$ws.Range("J66").Value = 4999.75
$ws.Range("K66").Value = 23432.4
$ws.Range("L66").Value = 44997.75
$ws.Range("M66").Value = -19688.4
$ws.Range("N66").Value = -52485.75

$ws.Range("H137").Value = 2657.9048
$ws.Range("J137").Value = 4138.8335
$ws.Range("L137").Value = 12416.5005
$ws.Range("N137").Value = -22616.5005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 11280.444
$ws.Range("I41").Value = 2575.5
$ws.Range("J41").Value = 18244.4
$ws.Range("K41").Value = 2575.5
$ws.Range("L41").Value = 18244.4
$ws.Range("M41").Value = -2220.5
$ws.Range("N41").Value = -18954.4

$ws.Range("H62").Value = 37940
$ws.Range("J62").Value = 37940
$ws.Range("L62").Value = 37940
$ws.Range("N62").Value = -39312

$ws.Range("H65").Value = 37940
$ws.Range("J65").Value = 37940
$ws.Range("L65").Value = 113820
$ws.Range("N65").Value = -120684

$ws.Range("H132").Value = 1892.0151
$ws.Range("I132").Value = 1243.4222
$ws.Range("J132").Value = 3281.8572
$ws.Range("K132").Value = 3730.2666
$ws.Range("L132").Value = 9845.5716
$ws.Range("M132").Value = -1200.2666
$ws.Range("N132").Value = -14905.5716

$ws.Range("H133").Value = 37872.777
$ws.Range("J133").Value = 37872.777
$ws.Range("L133").Value = 37872.777
$ws.Range("N133").Value = -47992.777

$ws.Range("H138").Value = 47803
$ws.Range("J138").Value = 47803
$ws.Range("L138").Value = 47803
$ws.Range("N138").Value = -58083

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 16002
$ws.Range("J26").Value = 29995
$ws.Range("L26").Value = 29995
$ws.Range("N26").Value = -30585

$ws.Range("H42").Value = 42498
$ws.Range("J42").Value = 42498
$ws.Range("L42").Value = 42498
$ws.Range("N42").Value = -43624

$ws.Range("H45").Value = 30942.5
$ws.Range("I45").Value = 21938.5
$ws.Range("J45").Value = 39946.5
$ws.Range("K45").Value = 21938.5
$ws.Range("L45").Value = 39946.5
$ws.Range("M45").Value = -21531.5
$ws.Range("N45").Value = -40760.5

$ws.Range("H49").Value = 42498
$ws.Range("J49").Value = 42498
$ws.Range("L49").Value = 42498
$ws.Range("N49").Value = -42792

$ws.Range("H76").Value = 22211.666
$ws.Range("J76").Value = 22211.666
$ws.Range("L76").Value = 22211.666
$ws.Range("N76").Value = -22887.666

$ws.Range("H79").Value = 22211.666
$ws.Range("J79").Value = 22211.666
$ws.Range("L79").Value = 22211.666
$ws.Range("N79").Value = -24551.666

$ws.Range("H132").Value = 8921.708
$ws.Range("I132").Value = 8958.032
$ws.Range("J132").Value = 8855.471
$ws.Range("K132").Value = 26874.096
$ws.Range("L132").Value = 26566.413
$ws.Range("M132").Value = -24344.096
$ws.Range("N132").Value = -31626.413

$ws.Range("H136").Value = 2261.125
$ws.Range("I136").Value = 1205.7778
$ws.Range("J136").Value = 7960
$ws.Range("K136").Value = 3617.3334
$ws.Range("L136").Value = 23880
$ws.Range("M136").Value = -1067.3334
$ws.Range("N136").Value = -28980

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3724.724
$ws.Range("I122").Value = 2517.5625
$ws.Range("J122").Value = 5210.4614
$ws.Range("K122").Value = 7552.6875
$ws.Range("L122").Value = 15631.3842
$ws.Range("M122").Value = -5102.6875
$ws.Range("N122").Value = -20531.3842

$ws.Range("H132").Value = 5378030
$ws.Range("I132").Value = 1186.7556
$ws.Range("K132").Value = 3560.2668
$ws.Range("M132").Value = -1030.2668
